$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns updated per row: E,F,G,H,I,J,M,N,O,P,Q,R,S,T
$cols = @("E","F","G","H","I","J","M","N","O","P","Q","R","S","T")

$data = @{
  2 = @(1, 0.3333333333333333, 0.022792, 0.06837600000000001, 0.001916327914826657, 0.001916327914826657, 10.34761366666667, 31.042841, 0.2299953477621856, 0.2299953477621856, 0.2358428106906667, 2.122585296216, 0.0004407465051969409, 0.000440746505196941)
  3 = @(1, 0.3333333333333333, 0.022792, 0.06837600000000001, 0.001916327914826657, 0.001916327914826657, 30.56986233333333, 91.709587, 0.6794731949692173, 0.6794731949692174, 0.6967483023013334, 6.270734720712, 0.001302093450895967, 0.001302093450895967)
  4 = @(1, 0.3333333333333333, 0.022792, 0.06837600000000001, 0.001916327914826657, 0.001916327914826657, 4.073058666666666, 12.219176, 0.09053145726859702, 0.09053145726859703, 0.09283315313066667, 0.835498378176, 0.0001734879587337491, 0.0001734879587337492)
  5 = @(3, 1, 4.402094666666667, 13.206284, 0.3701235913233977, 0.3701235913233977, 10.34761366666667, 31.042841, 0.2299953477621856, 0.2299953477621856, 45.55117493476044, 409.960574412844, 0.0851267041014139, 0.08512670410141393)
  6 = @(3, 1, 4.402094666666667, 13.206284, 0.3701235913233977, 0.3701235913233977, 30.56986233333333, 91.709587, 0.6794731949692173, 0.6794731949692174, 134.5714279383009, 1211.142851444708, 0.2514890591299899, 0.25148905912999)
  7 = @(3, 1, 4.402094666666667, 13.206284, 0.3701235913233977, 0.3701235913233977, 4.073058666666666, 12.219176, 0.09053145726859702, 0.09053145726859703, 17.92998983355378, 161.369908501984, 0.03350782809199385, 0.03350782809199385)
  8 = @(3, 1, 7.468693666666667, 22.406081, 0.6279600807617757, 0.6279600807617757, 10.34761366666667, 31.042841, 0.2299953477621856, 0.2299953477621856, 77.28315665734678, 695.548409916121, 0.1444278971555747, 0.1444278971555747)
  9 = @(3, 1, 7.468693666666667, 22.406081, 0.6279600807617757, 0.6279600807617757, 30.56986233333333, 91.709587, 0.6794731949692173, 0.6794731949692174, 228.3169371998385, 2054.852434798547, 0.4266820423883315, 0.4266820423883315)
  10 = @(3, 1, 7.468693666666667, 22.406081, 0.6279600807617757, 0.6279600807617757, 4.073058666666666, 12.219176, 0.09053145726859702, 0.09053145726859703, 30.42042746769511, 273.783847209256, 0.05685014121786942, 0.05685014121786943)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
